$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 111363020
$ws.Range("B2").Value2 = 78107
$ws.Range("D2").Value2 = "NT"
$ws.Range("E2").Value2 = 6453
$ws.Range("F2").Value2 = "Vedskivlav"
$ws.Range("G2").Value2 = "Hertelidea botryosa"
$ws.Range("H2").Value2 = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q2").Value2 = 593324.7367794912
$ws.Range("R2").Value2 = 6987171.102828567
$ws.Range("S2").Value2 = 10

# Row 4
$ws.Range("A4").Value2 = 111363022
$ws.Range("B4").Value2 = 77186
$ws.Range("D4").Value2 = "NT"
$ws.Range("E4").Value2 = 353
$ws.Range("F4").Value2 = "Dvärgbägarlav"
$ws.Range("G4").Value2 = "Cladonia parasitica"
$ws.Range("H4").Value2 = "(Hoffm.) Hoffm."
$ws.Range("Q4").Value2 = 593324.9051589288
$ws.Range("R4").Value2 = 6987181.108611984
$ws.Range("S4").Value2 = 10

# Row 5
$ws.Range("A5").Value2 = 111363030
$ws.Range("B5").Value2 = 77268
$ws.Range("D5").Value2 = "NT"
$ws.Range("E5").Value2 = 228912
$ws.Range("F5").Value2 = "Mörk kolflarnlav"
$ws.Range("G5").Value2 = "Carbonicola myrmecina"
$ws.Range("H5").Value2 = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q5").Value2 = 593355.1995546351
$ws.Range("R5").Value2 = 6987156.520171621
$ws.Range("S5").Value2 = 25

# Row 6
$ws.Range("A6").Value2 = 111363023
$ws.Range("B6").Value2 = 76918
$ws.Range("D6").Value2 = "NT"
$ws.Range("E6").Value2 = 6437
$ws.Range("F6").Value2 = "Blanksvart spiklav"
$ws.Range("G6").Value2 = "Calicium denigratum"
$ws.Range("H6").Value2 = "(Vain.) Tibell"
$ws.Range("Q6").Value2 = 593269.3631576585
$ws.Range("R6").Value2 = 6987149.513888635
$ws.Range("S6").Value2 = 10

# Row 7
$ws.Range("A7").Value2 = 111363026
$ws.Range("B7").Value2 = 90854
$ws.Range("D7").Value2 = "NT"
$ws.Range("E7").Value2 = 2079
$ws.Range("F7").Value2 = "Nordtagging"
$ws.Range("G7").Value2 = "Odonticium romellii"
$ws.Range("H7").Value2 = "(S.Lundell) Parmasto"
$ws.Range("Q7").Value2 = 593292.3890792141
$ws.Range("R7").Value2 = 6987203.815111163
$ws.Range("S7").Value2 = 10

# Row 8
$ws.Range("A8").Value2 = 111363028
$ws.Range("B8").Value2 = 77186
$ws.Range("D8").Value2 = "NT"
$ws.Range("E8").Value2 = 353
$ws.Range("F8").Value2 = "Dvärgbägarlav"
$ws.Range("G8").Value2 = "Cladonia parasitica"
$ws.Range("H8").Value2 = "(Hoffm.) Hoffm."
$ws.Range("Q8").Value2 = 593324.0129203054
$ws.Range("R8").Value2 = 6987101.07452714
$ws.Range("S8").Value2 = 10

# Row 9
$ws.Range("A9").Value2 = 111363021
$ws.Range("B9").Value2 = 89330
$ws.Range("D9").Value2 = "NT"
$ws.Range("E9").Value2 = 3242
$ws.Range("F9").Value2 = "Vitplätt"
$ws.Range("G9").Value2 = "Chaetodermella luna"
$ws.Range("H9").Value2 = "(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert"
$ws.Range("Q9").Value2 = 593278.356042281
$ws.Range("R9").Value2 = 6987153.408284122
$ws.Range("S9").Value2 = 10

# Row 10
$ws.Range("A10").Value2 = 111363025
$ws.Range("B10").Value2 = 89646
$ws.Range("D10").Value2 = "VU"
$ws.Range("E10").Value2 = 65
$ws.Range("F10").Value2 = "Fläckporing"
$ws.Range("G10").Value2 = "Anthoporia albobrunnea"
$ws.Range("H10").Value2 = "(Romell) Karasiński & Niemelä"
$ws.Range("Q10").Value2 = 593292.3890792141
$ws.Range("R10").Value2 = 6987203.815111163
$ws.Range("S10").Value2 = 10

# Row 11
$ws.Range("A11").Value2 = 111363029
$ws.Range("B11").Value2 = 76918
$ws.Range("D11").Value2 = "NT"
$ws.Range("E11").Value2 = 6437
$ws.Range("F11").Value2 = "Blanksvart spiklav"
$ws.Range("G11").Value2 = "Calicium denigratum"
$ws.Range("H11").Value2 = "(Vain.) Tibell"
$ws.Range("Q11").Value2 = 593312.9580448985
$ws.Range("R11").Value2 = 6987010.291132212
$ws.Range("S11").Value2 = 10

# Row 12
$ws.Range("A12").Value2 = 111363031
$ws.Range("B12").Value2 = 76918
$ws.Range("D12").Value2 = "NT"
$ws.Range("E12").Value2 = 6437
$ws.Range("F12").Value2 = "Blanksvart spiklav"
$ws.Range("G12").Value2 = "Calicium denigratum"
$ws.Range("H12").Value2 = "(Vain.) Tibell"
$ws.Range("Q12").Value2 = 593417.4633552339
$ws.Range("R12").Value2 = 6986985.556671137
$ws.Range("S12").Value2 = 10
